$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Dimension 1"
$ws.Range("C1").Value = "Dimension 2"
$ws.Range("D1").Value = "Dimension 3"
$ws.Range("E1").Value = "Dimension 4"

[void]$ws.Range("E1").Select()
